$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 489.44446
$ws.Range("I11").Value = 489.44446
$ws.Range("K11").Value = 489.44446
$ws.Range("M11").Value = -349.44446
$ws.Range("H70").Value = 989.4286
$ws.Range("J70").Value = 985
$ws.Range("L70").Value = 2955
$ws.Range("N70").Value = -3495
$ws.Range("H73").Value = 989.4286
$ws.Range("J73").Value = 985
$ws.Range("L73").Value = 2955
$ws.Range("N73").Value = -4827
$ws.Range("H74").Value = 6381.4443
$ws.Range("I74").Value = 6616.625
$ws.Range("K74").Value = 6616.625
$ws.Range("M74").Value = -5680.625
$ws.Range("H77").Value = 6381.4443
$ws.Range("I77").Value = 6616.625
$ws.Range("K77").Value = 33083.125
$ws.Range("M77").Value = -28403.125
$ws.Range("H132").Value = 1492.6129
$ws.Range("I132").Value = 1531.5186
$ws.Range("K132").Value = 4594.5558
$ws.Range("M132").Value = -2064.5558
$ws.Range("H137").Value = 2865364.8
$ws.Range("J137").Value = 7704591.5
$ws.Range("L137").Value = 23113774.5
$ws.Range("N137").Value = -23118874.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 210
$ws.Range("I4").Value = 310
$ws.Range("J4").Value = 110
$ws.Range("K4").Value = 310
$ws.Range("L4").Value = 110
$ws.Range("M4").Value = -194
$ws.Range("N4").Value = -342
$ws.Range("H32").Value = 349487.53
$ws.Range("I32").Value = 477934.66
$ws.Range("K32").Value = 477934.66
$ws.Range("M32").Value = -477647.66
$ws.Range("H43").Value = 36274.25
$ws.Range("J43").Value = 36274.25
$ws.Range("L43").Value = 36274.25
$ws.Range("N43").Value = -36900.25
$ws.Range("H61").Value = 1420551.4
$ws.Range("I61").Value = 47925.375
$ws.Range("J61").Value = 3154394.8
$ws.Range("K61").Value = 47925.375
$ws.Range("L61").Value = 3154394.8
$ws.Range("M61").Value = -47713.375
$ws.Range("N61").Value = -3154818.8
$ws.Range("H74").Value = 706790.75
$ws.Range("J74").Value = 1454707.8
$ws.Range("L74").Value = 1454707.8
$ws.Range("N74").Value = -1456455.8
$ws.Range("H77").Value = 706790.75
$ws.Range("J77").Value = 1454707.8
$ws.Range("L77").Value = 7273539
$ws.Range("N77").Value = -7282275
$ws.Range("H102").Value = 1926.6154
$ws.Range("I102").Value = 1831.5454
$ws.Range("K102").Value = 1831.5454
$ws.Range("M102").Value = -209.5454
$ws.Range("H122").Value = 1326
$ws.Range("I122").Value = 770
$ws.Range("J122").Value = 2252.6667
$ws.Range("K122").Value = 2310
$ws.Range("L122").Value = 6758.000100000001
$ws.Range("M122").Value = 140
$ws.Range("N122").Value = -11658.0001
$ws.Range("H126").Value = 9450
$ws.Range("I126").Value = 9450
$ws.Range("K126").Value = 28350
$ws.Range("M126").Value = -25880
$ws.Range("H136").Value = 1420551.4
$ws.Range("I136").Value = 47925.375
$ws.Range("J136").Value = 3154394.8
$ws.Range("K136").Value = 143776.125
$ws.Range("L136").Value = 9463184.399999999
$ws.Range("M136").Value = -141226.125
$ws.Range("N136").Value = -9468284.399999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6876.3076
$ws.Range("I86").Value = 5155.4546
$ws.Range("J86").Value = 8138.2666
$ws.Range("K86").Value = 5155.4546
$ws.Range("L86").Value = 8138.2666
$ws.Range("M86").Value = -4032.4546
$ws.Range("N86").Value = -10384.2666
$ws.Range("H89").Value = 6876.3076
$ws.Range("I89").Value = 5155.4546
$ws.Range("J89").Value = 8138.2666
$ws.Range("K89").Value = 25777.273
$ws.Range("L89").Value = 40691.333
$ws.Range("M89").Value = -20161.273
$ws.Range("N89").Value = -51923.333
$ws.Range("H95").Value = 38333.332
$ws.Range("J95").Value = 38333.332
$ws.Range("L95").Value = 38333.332
$ws.Range("N95").Value = -43825.332
$ws.Range("H99").Value = 13263.125
$ws.Range("I99").Value = 12420
$ws.Range("J99").Value = 16916.666
$ws.Range("K99").Value = 12420
$ws.Range("L99").Value = 16916.666
$ws.Range("M99").Value = -10922
$ws.Range("N99").Value = -19912.666
$ws.Range("H134").Value = 30002042
$ws.Range("I134").Value = 2243.0625
$ws.Range("K134").Value = 6729.1875
$ws.Range("M134").Value = -4194.1875

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3269.3704
$ws.Range("I31").Value = 4192
$ws.Range("K31").Value = 4192
$ws.Range("M31").Value = -3897
$ws.Range("H34").Value = 3269.3704
$ws.Range("I34").Value = 4192
$ws.Range("K34").Value = 4192
$ws.Range("M34").Value = -3990
$ws.Range("H132").Value = 1752.7142
$ws.Range("I132").Value = 1455.5555
$ws.Range("J132").Value = 2287.6
$ws.Range("K132").Value = 4366.666499999999
$ws.Range("L132").Value = 6862.799999999999
$ws.Range("M132").Value = -1836.666499999999
$ws.Range("N132").Value = -11922.8
$ws.Range("H134").Value = 2496.8333
$ws.Range("I134").Value = 2245.875
$ws.Range("K134").Value = 6737.625
$ws.Range("M134").Value = -4202.625

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 7582.3335
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 30000
$ws.Range("N32").Value = -30566
$ws.Range("H64").Value = 7534.6665
$ws.Range("I64").Value = 4861
$ws.Range("J64").Value = 9874.125
$ws.Range("K64").Value = 14583
$ws.Range("L64").Value = 29622.375
$ws.Range("M64").Value = -14313
$ws.Range("N64").Value = -30162.375
$ws.Range("H67").Value = 7534.6665
$ws.Range("I67").Value = 4861
$ws.Range("J67").Value = 9874.125
$ws.Range("K67").Value = 14583
$ws.Range("L67").Value = 29622.375
$ws.Range("M67").Value = -13647
$ws.Range("N67").Value = -31494.375
$ws.Range("H75").Value = 11907516
$ws.Range("I75").Value = 832
$ws.Range("J75").Value = 17860858
$ws.Range("K75").Value = 2496
$ws.Range("L75").Value = 53582574
$ws.Range("M75").Value = -1498
$ws.Range("N75").Value = -53584570
$ws.Range("H78").Value = 11907516
$ws.Range("I78").Value = 832
$ws.Range("J78").Value = 17860858
$ws.Range("K78").Value = 7488
$ws.Range("L78").Value = 160747722
$ws.Range("M78").Value = -2496
$ws.Range("N78").Value = -160757706

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.25
$ws.Range("I2").Value = 81.42856999999999
$ws.Range("J2").Value = 100.2
$ws.Range("K2").Value = 81.42856999999999
$ws.Range("L2").Value = 100.2
$ws.Range("M2").Value = 31.57143000000001
$ws.Range("N2").Value = -326.2
$ws.Range("H86").Value = 98000
$ws.Range("J86").Value = 98000
$ws.Range("L86").Value = 98000
$ws.Range("N86").Value = -100372
$ws.Range("H89").Value = 98000
$ws.Range("J89").Value = 98000
$ws.Range("L89").Value = 294000
$ws.Range("N89").Value = -305856
$ws.Range("H126").Value = 3658.75
$ws.Range("I126").Value = 1540.3334
$ws.Range("K126").Value = 4621.0002
$ws.Range("M126").Value = -2151.0002

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3206.6667
$ws.Range("I22").Value = 723.3333
$ws.Range("K22").Value = 723.3333
$ws.Range("M22").Value = -428.3333
$ws.Range("H27").Value = 3206.6667
$ws.Range("I27").Value = 723.3333
$ws.Range("K27").Value = 723.3333
$ws.Range("M27").Value = -616.3333
$ws.Range("H46").Value = 7999.0454
$ws.Range("J46").Value = 3577.9473
$ws.Range("L46").Value = 3577.9473
$ws.Range("N46").Value = -3953.9473
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 73866.92999999999
$ws.Range("I136").Value = 113596.445
$ws.Range("J136").Value = 2353.8
$ws.Range("K136").Value = 340789.335
$ws.Range("L136").Value = 7061.400000000001
$ws.Range("M136").Value = -338239.335
$ws.Range("N136").Value = -12161.4

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 83583.336
$ws.Range("J46").Value = 83583.336
$ws.Range("L46").Value = 83583.336
$ws.Range("N46").Value = -84045.336
$ws.Range("H62").Value = 3167.6667
$ws.Range("I62").Value = 3001.5
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3001.5
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2377.5
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3167.6667
$ws.Range("I65").Value = 3001.5
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 15007.5
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -11887.5
$ws.Range("N65").Value = -23740
$ws.Range("H95").Value = 51999.668
$ws.Range("J95").Value = 51999.668
$ws.Range("L95").Value = 51999.668
$ws.Range("N95").Value = -57491.668
$ws.Range("H122").Value = 1277.4231
$ws.Range("I122").Value = 1060.5
$ws.Range("J122").Value = 2000.5
$ws.Range("K122").Value = 3181.5
$ws.Range("L122").Value = 6001.5
$ws.Range("M122").Value = -731.5
$ws.Range("N122").Value = -10901.5
$ws.Range("H134").Value = 83583.336
$ws.Range("J134").Value = 83583.336
$ws.Range("L134").Value = 250750.008
$ws.Range("N134").Value = -255820.008
$ws.Range("H136").Value = 665.63635
$ws.Range("J136").Value = 1037.2
$ws.Range("L136").Value = 3111.6
$ws.Range("N136").Value = -8211.6
